# Commit: "Finished experiment 3 and tryid making agent judge more robust"
# Mark experiments 2 and 3 (rows 3 and 4) as Done = TRUE in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 -> Experiment 2, Row 4 -> Experiment 3 : set "Done" column (G) to TRUE
$ws.Range("G3").Value = $true
$ws.Range("G4").Value = $true

# Update the active cell selection to reflect where the user left off editing
$ws.Range("H7").Select()
